# Update the heading date
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-05-26 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-05-27 Tuesday", 2)

# Replace each arithmetic-problem cell in the 20x5 table with the new value,
# preserving all cell/run formatting (font, size, paragraph alignment, etc.)
$newValues = @(
    @("40+18=", "93-31=", "99-5=", "77-33=", "88-3="),
    @("48+41=", "28+62=", "62-57=", "47-30=", "77-53="),
    @("29-11=", "85-69=", "27-0=", "81-81=", "60+20="),
    @("88-34=", "43-26=", "93-1=", "64-12=", "10+2="),
    @("56-6=", "64-25=", "61+34=", "71+4=", "38-19="),
    @("46-1=", "53-25=", "43+9=", "67-11=", "31+28="),
    @("34+39=", "11+48=", "77-51=", "70+24=", "71-4="),
    @("24+15=", "21+54=", "52-18=", "1+84=", "16+28="),
    @("56+5=", "23-4=", "19-18=", "19+25=", "10+74="),
    @("58-47=", "40+29=", "61+13=", "60-38=", "25+27="),
    @("15+61=", "23-10=", "75+10=", "47+32=", "42+23="),
    @("27+66=", "83-16=", "81-14=", "43-43=", "60-6="),
    @("4+24=", "45+41=", "83-12=", "59-39=", "34+40="),
    @("3+1=", "31+41=", "0+23=", "32+27=", "41-8="),
    @("80-64=", "32+9=", "90-42=", "65+23=", "56+31="),
    @("13+70=", "84-9=", "0+54=", "3+11=", "10+33="),
    @("76-24=", "22+7=", "50-9=", "50+40=", "76-26="),
    @("45+32=", "50+46=", "67-40=", "52-43=", "49+2="),
    @("67-24=", "52-33=", "36-18=", "85-79=", "20+63="),
    @("31+30=", "0+34=", "39-24=", "28+67=", "66-16=")
)

$t = $d.Tables(1)
for ($r = 1; $r -le $newValues.Length; $r++) {
    $row = $newValues[$r - 1]
    for ($c = 1; $c -le $row.Length; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $row[$c - 1]
    }
}
